# "Scrape new links of Ilca 6 & Ilca 7" - populate the Flotilha (fleet)
# column (F) for every result row: rows that raced in the Gold fleet
# (previously blank) are now tagged "G", and the rows that were flagged
# "medal" (medal race) are now tagged "MR".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp

for ($r = 2; $r -le $lastRow; $r++) {
    $flotilhaCell = $ws.Cells.Item($r, 6)  # column F
    if ($flotilhaCell.Value() -eq "medal") {
        $flotilhaCell.Value = "MR"
    } else {
        $flotilhaCell.Value = "G"
    }
}
